# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (same per-fund layout as the other
#    quarter sheets) right before the "总计" (totals) summary sheet.
# 2) Prepend a "2022-Q1" row to the "总计" sheet, shifting the previous
#    rows down by one and renumbering the index column (A).

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) New "2022-Q1" worksheet, inserted immediately before "总计"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Clone the formatting (fonts/borders/alignment) of an existing, already
# structured quarter sheet so the new sheet's styles line up with the rest
# of the workbook instead of minting brand-new style records.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H3").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "290012"
$newSheet.Range("C2").Value = "泰信行业精选灵活配置混合A"
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.76"
$newSheet.Range("E2").Value = "92.62"
$newSheet.Range("F2").Value = "5.15"
$newSheet.Range("G2").Value = "0.0391"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "002583"
$newSheet.Range("C3").Value = "泰信行业精选灵活配置混合C"
$newSheet.Range("D3:F3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.00"
$newSheet.Range("E3").Value = "92.62"
$newSheet.Range("F3").Value = "5.15"
$newSheet.Range("G3").Value = 0
$newSheet.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2) "总计" sheet: push existing quarter rows down one row and insert the
#    new 2022-Q1 totals at the top of the data (row 2).
#
# NB: re-resolve the "总计" worksheet by name now that the sheet
# collection has changed shape -- the handle captured before
# Worksheets.Add() above tracks the *slot*, not the sheet, and that slot
# now belongs to the freshly inserted "2022-Q1" sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$lastRow = 6  # existing data occupies rows 2..6 (2021-Q4 .. 2020-Q4)

$dates = @()
$counts = @()
$values = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $dates += $total.Cells.Item($r, 2).Value2
    $counts += $total.Cells.Item($r, 3).Value2
    $values += $total.Cells.Item($r, 4).Value2
}

for ($i = 0; $i -lt $dates.Count; $i++) {
    $targetRow = $i + 3
    $total.Cells.Item($targetRow, 1).Value = $i + 1
    $total.Cells.Item($targetRow, 2).Value = $dates[$i]
    $total.Cells.Item($targetRow, 3).Value = $counts[$i]
    $total.Cells.Item($targetRow, 4).Value = $values[$i]
}

# Carry the bold/bordered index-column style onto the newly created last
# row (it falls outside the sheet's original used range, so it wouldn't
# otherwise pick up the "A" column formatting).
$total.Range("A" + $lastRow).Copy()
$total.Range("A" + ($lastRow + 1)).PasteSpecial(-4122)
$total.Cells.Item($lastRow + 1, 1).Value = $dates.Count

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.04

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
